# Apply row-content reshuffle to rows 2-22 of the "Artfynd" sheet.
# Each row keeps its constant columns (C, P, S, T, U, V, W, AD, AE, AG, AW, AX, etc.)
# but receives a new set of record values (A, B, D, E, F, G, H, Q, R, Y, Z, AA, AB)
# taken from another row of the original data (a full re-sort/re-sync of the report).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111808387
$ws.Range("B2").Value = 89419
$ws.Range("E2").Value = 1204
$ws.Range("Q2").Value = 611016.8359391808
$ws.Range("R2").Value = 7181141.984797659
$ws.Range("D2").Value = '''NT'
$ws.Range("F2").Value = '''Gränsticka'
$ws.Range("G2").Value = '''Phellopilus nigrolimitatus'
$ws.Range("H2").Value = '''(Romell) Niemelä, T.Wagner & M.Fisch.'
$ws.Range("Y2").Value = '''2023-08-31'
$ws.Range("Z2").Value = '''13:26'
$ws.Range("AA2").Value = '''2023-08-31'
$ws.Range("AB2").Value = '''13:26'

# Row 3
$ws.Range("A3").Value = 111808817
$ws.Range("B3").Value = 56398
$ws.Range("E3").Value = 100109
$ws.Range("Q3").Value = 610921.7319367616
$ws.Range("R3").Value = 7180934.079081071
$ws.Range("D3").Value = '''NT'
$ws.Range("F3").Value = '''Tretåig hackspett'
$ws.Range("G3").Value = '''Picoides tridactylus'
$ws.Range("H3").Value = '''(Linnaeus, 1758)'
$ws.Range("Y3").Value = '''2023-08-31'
$ws.Range("Z3").Value = '''13:45'
$ws.Range("AA3").Value = '''2023-08-31'
$ws.Range("AB3").Value = '''13:45'

# Row 4
$ws.Range("A4").Value = 111809897
$ws.Range("B4").Value = 85715
$ws.Range("E4").Value = 510
$ws.Range("Q4").Value = 610718.603132805
$ws.Range("R4").Value = 7180857.334717941
$ws.Range("D4").Value = '''NT'
$ws.Range("F4").Value = '''Doftskinn'
$ws.Range("G4").Value = '''Cystostereum murrayi'
$ws.Range("H4").Value = '''(Berk. & M.A. Curtis.) Pouzar'
$ws.Range("Y4").Value = '''2023-08-31'
$ws.Range("Z4").Value = '''14:50'
$ws.Range("AA4").Value = '''2023-08-31'
$ws.Range("AB4").Value = '''14:50'

# Row 5
$ws.Range("A5").Value = 111806969
$ws.Range("B5").Value = 56398
$ws.Range("E5").Value = 100109
$ws.Range("Q5").Value = 610695.5210812307
$ws.Range("R5").Value = 7181007.871029559
$ws.Range("D5").Value = '''NT'
$ws.Range("F5").Value = '''Tretåig hackspett'
$ws.Range("G5").Value = '''Picoides tridactylus'
$ws.Range("H5").Value = '''(Linnaeus, 1758)'
$ws.Range("Y5").Value = '''2023-08-31'
$ws.Range("Z5").Value = '''12:09'
$ws.Range("AA5").Value = '''2023-08-31'
$ws.Range("AB5").Value = '''12:09'

# Row 6
$ws.Range("A6").Value = 111808328
$ws.Range("B6").Value = 90087
$ws.Range("E6").Value = 3298
$ws.Range("Q6").Value = 611008.4619706698
$ws.Range("R6").Value = 7181032.205813259
$ws.Range("D6").Value = '''LC'
$ws.Range("F6").Value = '''Trådticka'
$ws.Range("G6").Value = '''Climacocystis borealis'
$ws.Range("H6").Value = '''(Fr.) Kotl. & Pouzar'
$ws.Range("Y6").Value = '''2023-08-31'
$ws.Range("Z6").Value = '''13:23'
$ws.Range("AA6").Value = '''2023-08-31'
$ws.Range("AB6").Value = '''13:23'

# Row 7
$ws.Range("A7").Value = 111807777
$ws.Range("B7").Value = 56398
$ws.Range("E7").Value = 100109
$ws.Range("Q7").Value = 610955.5779051523
$ws.Range("R7").Value = 7181003.318500374
$ws.Range("D7").Value = '''NT'
$ws.Range("F7").Value = '''Tretåig hackspett'
$ws.Range("G7").Value = '''Picoides tridactylus'
$ws.Range("H7").Value = '''(Linnaeus, 1758)'
$ws.Range("Y7").Value = '''2023-08-31'
$ws.Range("Z7").Value = '''12:48'
$ws.Range("AA7").Value = '''2023-08-31'
$ws.Range("AB7").Value = '''12:48'

# Row 8
$ws.Range("A8").Value = 111809318
$ws.Range("B8").Value = 56398
$ws.Range("E8").Value = 100109
$ws.Range("Q8").Value = 610793.1422167002
$ws.Range("R8").Value = 7180891.290986502
$ws.Range("D8").Value = '''NT'
$ws.Range("F8").Value = '''Tretåig hackspett'
$ws.Range("G8").Value = '''Picoides tridactylus'
$ws.Range("H8").Value = '''(Linnaeus, 1758)'
$ws.Range("Y8").Value = '''2023-08-31'
$ws.Range("Z8").Value = '''14:15'
$ws.Range("AA8").Value = '''2023-08-31'
$ws.Range("AB8").Value = '''14:15'

# Row 9
$ws.Range("A9").Value = 111808022
$ws.Range("B9").Value = 90087
$ws.Range("E9").Value = 3298
$ws.Range("Q9").Value = 611005.4375152331
$ws.Range("R9").Value = 7181032.949711116
$ws.Range("D9").Value = '''LC'
$ws.Range("F9").Value = '''Trådticka'
$ws.Range("G9").Value = '''Climacocystis borealis'
$ws.Range("H9").Value = '''(Fr.) Kotl. & Pouzar'
$ws.Range("Y9").Value = '''2023-08-31'
$ws.Range("Z9").Value = '''12:58'
$ws.Range("AA9").Value = '''2023-08-31'
$ws.Range("AB9").Value = '''12:58'

# Row 10
$ws.Range("A10").Value = 111807814
$ws.Range("B10").Value = 56398
$ws.Range("E10").Value = 100109
$ws.Range("Q10").Value = 610955.7664983921
$ws.Range("R10").Value = 7180998.193866283
$ws.Range("D10").Value = '''NT'
$ws.Range("F10").Value = '''Tretåig hackspett'
$ws.Range("G10").Value = '''Picoides tridactylus'
$ws.Range("H10").Value = '''(Linnaeus, 1758)'
$ws.Range("Y10").Value = '''2023-08-31'
$ws.Range("Z10").Value = '''12:50'
$ws.Range("AA10").Value = '''2023-08-31'
$ws.Range("AB10").Value = '''12:50'

# Row 11
$ws.Range("A11").Value = 111809026
$ws.Range("B11").Value = 56398
$ws.Range("E11").Value = 100109
$ws.Range("Q11").Value = 610772.3941544103
$ws.Range("R11").Value = 7180884.969268824
$ws.Range("D11").Value = '''NT'
$ws.Range("F11").Value = '''Tretåig hackspett'
$ws.Range("G11").Value = '''Picoides tridactylus'
$ws.Range("H11").Value = '''(Linnaeus, 1758)'
$ws.Range("Y11").Value = '''2023-08-31'
$ws.Range("Z11").Value = '''13:56'
$ws.Range("AA11").Value = '''2023-08-31'
$ws.Range("AB11").Value = '''13:56'

# Row 12
$ws.Range("A12").Value = 111808515
$ws.Range("B12").Value = 56398
$ws.Range("E12").Value = 100109
$ws.Range("Q12").Value = 610986.3631281323
$ws.Range("R12").Value = 7181120.765008576
$ws.Range("D12").Value = '''NT'
$ws.Range("F12").Value = '''Tretåig hackspett'
$ws.Range("G12").Value = '''Picoides tridactylus'
$ws.Range("H12").Value = '''(Linnaeus, 1758)'
$ws.Range("Y12").Value = '''2023-08-31'
$ws.Range("Z12").Value = '''13:29'
$ws.Range("AA12").Value = '''2023-08-31'
$ws.Range("AB12").Value = '''13:29'

# Row 13
$ws.Range("A13").Value = 111808957
$ws.Range("B13").Value = 56398
$ws.Range("E13").Value = 100109
$ws.Range("Q13").Value = 610787.509024983
$ws.Range("R13").Value = 7180858.155172224
$ws.Range("D13").Value = '''NT'
$ws.Range("F13").Value = '''Tretåig hackspett'
$ws.Range("G13").Value = '''Picoides tridactylus'
$ws.Range("H13").Value = '''(Linnaeus, 1758)'
$ws.Range("Y13").Value = '''2023-08-31'
$ws.Range("Z13").Value = '''13:52'
$ws.Range("AA13").Value = '''2023-08-31'
$ws.Range("AB13").Value = '''13:52'

# Row 14
$ws.Range("A14").Value = 111807370
$ws.Range("B14").Value = 56398
$ws.Range("E14").Value = 100109
$ws.Range("Q14").Value = 610667.3767981895
$ws.Range("R14").Value = 7181039.764941735
$ws.Range("D14").Value = '''NT'
$ws.Range("F14").Value = '''Tretåig hackspett'
$ws.Range("G14").Value = '''Picoides tridactylus'
$ws.Range("H14").Value = '''(Linnaeus, 1758)'
$ws.Range("Y14").Value = '''2023-08-31'
$ws.Range("Z14").Value = '''12:36'
$ws.Range("AA14").Value = '''2023-08-31'
$ws.Range("AB14").Value = '''12:36'

# Row 15
$ws.Range("A15").Value = 111807821
$ws.Range("B15").Value = 56398
$ws.Range("E15").Value = 100109
$ws.Range("Q15").Value = 610967.5972068857
$ws.Range("R15").Value = 7181002.477957427
$ws.Range("D15").Value = '''NT'
$ws.Range("F15").Value = '''Tretåig hackspett'
$ws.Range("G15").Value = '''Picoides tridactylus'
$ws.Range("H15").Value = '''(Linnaeus, 1758)'
$ws.Range("Y15").Value = '''2023-08-31'
$ws.Range("Z15").Value = '''12:51'
$ws.Range("AA15").Value = '''2023-08-31'
$ws.Range("AB15").Value = '''12:51'

# Row 16
$ws.Range("A16").Value = 111807055
$ws.Range("B16").Value = 77515
$ws.Range("E16").Value = 6425
$ws.Range("Q16").Value = 610666.4119294117
$ws.Range("R16").Value = 7181042.722880279
$ws.Range("D16").Value = '''NT'
$ws.Range("F16").Value = '''Garnlav'
$ws.Range("G16").Value = '''Alectoria sarmentosa'
$ws.Range("H16").Value = '''(Ach.) Ach.'
$ws.Range("Y16").Value = '''2023-08-31'
$ws.Range("Z16").Value = '''12:10'
$ws.Range("AA16").Value = '''2023-08-31'
$ws.Range("AB16").Value = '''12:10'

# Row 17
$ws.Range("A17").Value = 111807092
$ws.Range("B17").Value = 77515
$ws.Range("E17").Value = 6425
$ws.Range("Q17").Value = 610678.3995887624
$ws.Range("R17").Value = 7181042.735301252
$ws.Range("D17").Value = '''NT'
$ws.Range("F17").Value = '''Garnlav'
$ws.Range("G17").Value = '''Alectoria sarmentosa'
$ws.Range("H17").Value = '''(Ach.) Ach.'
$ws.Range("Y17").Value = '''2023-08-31'
$ws.Range("Z17").Value = '''12:11'
$ws.Range("AA17").Value = '''2023-08-31'
$ws.Range("AB17").Value = '''12:11'

# Row 18
$ws.Range("A18").Value = 111808659
$ws.Range("B18").Value = 56398
$ws.Range("E18").Value = 100109
$ws.Range("Q18").Value = 610912.6426496292
$ws.Range("R18").Value = 7180936.738188162
$ws.Range("D18").Value = '''NT'
$ws.Range("F18").Value = '''Tretåig hackspett'
$ws.Range("G18").Value = '''Picoides tridactylus'
$ws.Range("H18").Value = '''(Linnaeus, 1758)'
$ws.Range("Y18").Value = '''2023-08-31'
$ws.Range("Z18").Value = '''13:39'
$ws.Range("AA18").Value = '''2023-08-31'
$ws.Range("AB18").Value = '''13:39'

# Row 19
$ws.Range("A19").Value = 111808000
$ws.Range("B19").Value = 56398
$ws.Range("E19").Value = 100109
$ws.Range("Q19").Value = 610972.8318012832
$ws.Range("R19").Value = 7180976.585010829
$ws.Range("D19").Value = '''NT'
$ws.Range("F19").Value = '''Tretåig hackspett'
$ws.Range("G19").Value = '''Picoides tridactylus'
$ws.Range("H19").Value = '''(Linnaeus, 1758)'
$ws.Range("Y19").Value = '''2023-08-31'
$ws.Range("Z19").Value = '''12:56'
$ws.Range("AA19").Value = '''2023-08-31'
$ws.Range("AB19").Value = '''12:56'

# Row 20
$ws.Range("A20").Value = 111808532
$ws.Range("B20").Value = 89423
$ws.Range("E20").Value = 5432
$ws.Range("Q20").Value = 610963.6803355663
$ws.Range("R20").Value = 7181097.266412931
$ws.Range("D20").Value = '''NT'
$ws.Range("F20").Value = '''Granticka'
$ws.Range("G20").Value = '''Porodaedalea chrysoloma'
$ws.Range("H20").Value = '''(Fr.) Fiasson & Niemelä'
$ws.Range("Y20").Value = '''2023-08-31'
$ws.Range("Z20").Value = '''13:32'
$ws.Range("AA20").Value = '''2023-08-31'
$ws.Range("AB20").Value = '''13:32'

# Row 21
$ws.Range("A21").Value = 111808676
$ws.Range("B21").Value = 56398
$ws.Range("E21").Value = 100109
$ws.Range("Q21").Value = 610906.3771198173
$ws.Range("R21").Value = 7180955.75143602
$ws.Range("D21").Value = '''NT'
$ws.Range("F21").Value = '''Tretåig hackspett'
$ws.Range("G21").Value = '''Picoides tridactylus'
$ws.Range("H21").Value = '''(Linnaeus, 1758)'
$ws.Range("Y21").Value = '''2023-08-31'
$ws.Range("Z21").Value = '''13:40'
$ws.Range("AA21").Value = '''2023-08-31'
$ws.Range("AB21").Value = '''13:40'

# Row 22
$ws.Range("A22").Value = 111806888
$ws.Range("B22").Value = 56398
$ws.Range("E22").Value = 100109
$ws.Range("Q22").Value = 610667.2784075569
$ws.Range("R22").Value = 7181030.781246905
$ws.Range("D22").Value = '''NT'
$ws.Range("F22").Value = '''Tretåig hackspett'
$ws.Range("G22").Value = '''Picoides tridactylus'
$ws.Range("H22").Value = '''(Linnaeus, 1758)'
$ws.Range("Y22").Value = '''1970-01-01'
$ws.Range("Z22").Value = '''00:00'
$ws.Range("AA22").Value = '''1970-01-01'
$ws.Range("AB22").Value = '''00:00'

